# Update the "Förändrad" (Changed) date column C for rows 2-14
# from 45224 (2023-10-25) to 45233 (2023-11-03), per the commit's
# automatic update of files.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45224) {
        $cell.Value2 = 45233
    }
}
